{"js": "// Replace the 25 \"NNN\u00d7N=\" multiplication prompts in the practice-sheet\n// table with the new values from the authoring commit. Every old value is\n// a unique string in the document, so a plain exact-text search/replace\n// for each pair is safe and unambiguous.\nconst replacements = [\n  [\"114\u00d73=\", \"878\u00d74=\"],\n  [\"685\u00d76=\", \"264\u00d79=\"],\n  [\"106\u00d79=\", \"297\u00d79=\"],\n  [\"530\u00d77=\", \"358\u00d74=\"],\n  [\"744\u00d76=\", \"380\u00d78=\"],\n  [\"661\u00d75=\", \"646\u00d77=\"],\n  [\"680\u00d78=\", \"628\u00d74=\"],\n  [\"494\u00d75=\", \"839\u00d74=\"],\n  [\"127\u00d72=\", \"103\u00d74=\"],\n  [\"315\u00d76=\", \"999\u00d79=\"],\n  [\"278\u00d73=\", \"342\u00d77=\"],\n  [\"800\u00d76=\", \"847\u00d79=\"],\n  [\"736\u00d77=\", \"802\u00d75=\"],\n  [\"197\u00d72=\", \"562\u00d77=\"],\n  [\"447\u00d77=\", \"740\u00d74=\"],\n  [\"751\u00d73=\", \"345\u00d75=\"],\n  [\"799\u00d75=\", \"573\u00d77=\"],\n  [\"675\u00d72=\", \"443\u00d77=\"],\n  [\"782\u00d78=\", \"208\u00d77=\"],\n  [\"455\u00d77=\", \"120\u00d75=\"],\n  [\"503\u00d78=\", \"718\u00d73=\"],\n  [\"436\u00d73=\", \"721\u00d79=\"],\n  [\"365\u00d78=\", \"418\u00d73=\"],\n  [\"760\u00d78=\", \"546\u00d78=\"],\n  [\"899\u00d73=\", \"587\u00d74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"NNN\u00d7N=\" multiplication prompts in the practice-sheet\n# table with the new values from the authoring commit. Every old value is\n# a unique string in the document, so Find/Replace (wdReplaceAll, scoped\n# to the whole document) for each pair is safe and unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"114\u00d73=\", \"878\u00d74=\"),\n    @(\"685\u00d76=\", \"264\u00d79=\"),\n    @(\"106\u00d79=\", \"297\u00d79=\"),\n    @(\"530\u00d77=\", \"358\u00d74=\"),\n    @(\"744\u00d76=\", \"380\u00d78=\"),\n    @(\"661\u00d75=\", \"646\u00d77=\"),\n    @(\"680\u00d78=\", \"628\u00d74=\"),\n    @(\"494\u00d75=\", \"839\u00d74=\"),\n    @(\"127\u00d72=\", \"103\u00d74=\"),\n    @(\"315\u00d76=\", \"999\u00d79=\"),\n    @(\"278\u00d73=\", \"342\u00d77=\"),\n    @(\"800\u00d76=\", \"847\u00d79=\"),\n    @(\"736\u00d77=\", \"802\u00d75=\"),\n    @(\"197\u00d72=\", \"562\u00d77=\"),\n    @(\"447\u00d77=\", \"740\u00d74=\"),\n    @(\"751\u00d73=\", \"345\u00d75=\"),\n    @(\"799\u00d75=\", \"573\u00d77=\"),\n    @(\"675\u00d72=\", \"443\u00d77=\"),\n    @(\"782\u00d78=\", \"208\u00d77=\"),\n    @(\"455\u00d77=\", \"120\u00d75=\"),\n    @(\"503\u00d78=\", \"718\u00d73=\"),\n    @(\"436\u00d73=\", \"721\u00d79=\"),\n    @(\"365\u00d78=\", \"418\u00d73=\"),\n    @(\"760\u00d78=\", \"546\u00d78=\"),\n    @(\"899\u00d73=\", \"587\u00d74=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
